$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("welcome")

# --- 1. Fix up per-row cell styles (A:B) BEFORE touching values / deleting rows,
#     using rows whose style already matches what the target needs as templates.
#     Row 2 = style set (A=6,B=5); Row 3 = style set (A=2,B=4) - both keep the
#     same style in the final layout, so they are safe, stable templates.
#     Row 14/15 carry the (A=7,B=4) / (A=8,B=5) style sets still needed by the
#     final rows 12/13 before those source rows get deleted.
$ws.Range("A2:B2").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)

$ws.Range("A3:B3").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$ws.Range("A3:B3").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A14:B14").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A15:B15").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 2. Drop the two trailing rows (14 & 15) - final sheet only has 13 rows.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# --- 3. Clear the now-unwanted C9/C10 cells (final rows 9-13 have no C cell).
$ws.Range("C9").Clear()
$ws.Range("C10").Clear()

# --- 4. Write the final text content for every row.
$ws.Range("B1").Value = "陌生人"

$ws.Range("A2").Value = "DESCRIPTIONTEXT"
$ws.Range("B2").Value = "你正在读易卜生的名著《海达·高布乐》，一个陌生人在你旁边"

$ws.Range("A3").Value = "WELCOMEWORD"
$ws.Range("B3").Value = "你好，陌生人，你看过易卜生的《海达·高布乐》吗？"

$ws.Range("A4").Value = "NO_FOCUS_TOPIC"
$ws.Range("B4").Value = "对陌生人的话并不感兴趣，你劝陌生人也去看看《海达·高布乐》这部剧，你"

$ws.Range("A5").Value = "challenge"
$ws.Range("B5").Value = "SOLID哎……自以为是的人类啊[困]"

$ws.Range("A6").Value = "greeting"
$ws.Range("B6").Value = "想问问他对于《海达·高布乐》这本书的看法，你"

$ws.Range("A7").Value = "bye"
$ws.Range("B7").Value = "SOLID嗯，那我继续看书了，祝你开心[调皮]"

$ws.Range("A8").Value = "noknowledge"
$ws.Range("B8").Value = "想劝陌生人先去看看《海达·高布乐》然后再跟你对话，你"

$ws.Range("A9").Value = "notinterested"
$ws.Range("B9").Value = "SOLID好吧，那先这样，我们下次再聊[开心]"

$ws.Range("A10").Value = "talkabout"
$ws.Range("B10").Value = "对陌生人的话很感兴趣，追问"

$ws.Range("A11").Value = "interestedin"
$ws.Range("B11").Value = "继续"

$ws.Range("A12").Value = "misunderstand"
$ws.Range("B12").Value = "发现自己刚才说错了，于是赶忙遮掩"

$ws.Range("A13").Value = "challenge_bye"
$ws.Range("B13").Value = "SOLID我承认AI现在还是发展期，但你这样说还是让我很伤心，不想再理你了[委屈]"

# --- 5. Row heights for the final layout.
$ws.Rows.Item(7).RowHeight = 25
$ws.Rows.Item(8).RowHeight = 24
$ws.Rows.Item(9).RowHeight = 25
$ws.Rows.Item(11).RowHeight = 24
$ws.Rows.Item(13).RowHeight = 25

# --- 6. Selection moves to B22 on the welcome sheet.
$ws.Range("B22").Select()
